$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the <b>…</b> emphasis markup that used to wrap "stand/light" keyword
# hits inside every listing title (column B).
$ws.Cells.Item(2, 2).Value = "캥거 북유럽 플로어 스탠드 조명 장스탠드 LED 거실등 KG-019"
$ws.Cells.Item(3, 2).Value = "파파조명 다이아 장스탠드 거실 스탠드 조명 무드등"
$ws.Cells.Item(4, 2).Value = "이케아 스탠드조명 장스탠드 플로어 독서등 거실 토가르프 레르스타 헥토그람 오르스티드 알렝"
$ws.Cells.Item(5, 2).Value = "이케아 스탠드 조명 LED거실등 무드등 장스탠드 침실 인테리어"
$ws.Cells.Item(6, 2).Value = "이케아 TAGARP 토가르프 상향식플로어스탠드 조명 스텐드 전구미포함"
$ws.Cells.Item(7, 2).Value = "예쁜 거실 인테리어 스탠드조명 스탠드등 장스탠드 무선리모컨"
$ws.Cells.Item(8, 2).Value = "이케아 스탠드 조명 등 LERSTA"
$ws.Cells.Item(9, 2).Value = "이케아 ARSTID 오르스티드 무드등 플로어스탠드 장스탠드 거실조명 엔틱"
$ws.Cells.Item(10, 2).Value = "무선 LED 장스탠드 플로어 롱 거실 침대 스탠드 조명 무드등 램프"
$ws.Cells.Item(11, 2).Value = "장스탠드 인테리어스탠드 스탠드조명 플로어램프 스탠드조명"
$ws.Cells.Item(12, 2).Value = "안방 조명 장 스탠드 북유럽 모던 거실 침실 침대 취침 플로어 롱 깃털 무드 일자형 LED조명"
$ws.Cells.Item(13, 2).Value = "이케아 IKEA 장스탠드 LED 스탠드 조명 독서등 무드등 거실등"
$ws.Cells.Item(14, 2).Value = "장스탠드 무드등 레르스타 뉘모네 조명등 독서등 플로어스탠드 거실장 이케아 활장 스탠딩조명"
$ws.Cells.Item(15, 2).Value = "RUSTA 마켓비 장스탠드 거실조명 스탠딩"
$ws.Cells.Item(16, 2).Value = "이케아 조명 플로어스탠드 LERSTA 거실 202.842.48"
$ws.Cells.Item(17, 2).Value = "북유럽 원목 장 스탠드 조명 선반 거실 침실 조명"
$ws.Cells.Item(18, 2).Value = "스탠드조명 장스탠드 램프 이케아 거실 조명등 무드등 플로어 스탠딩"
$ws.Cells.Item(19, 2).Value = "PH5 플로어스탠드 침실 거실 무드등 장스탠드 조명 카페 장조명 가로등"
$ws.Cells.Item(20, 2).Value = "파인굿즈 LED 장스탠드 조명 거실 인테리어 책상 학습 독서 스탠드조명"
$ws.Cells.Item(21, 2).Value = "북유럽 탁자 장 스탠드조명 거실 침실 테이블"
$ws.Cells.Item(22, 2).Value = "거실 스탠드조명 심플 모던 북유럽 크리에이티브 가정용 침실 헤드라이트 개성"
$ws.Cells.Item(23, 2).Value = "엘사 크리스탈 눈꽃스탠드 조명 램프11 이벤트행사"
$ws.Cells.Item(24, 2).Value = "프리미엄 깃털 장스탠드 인테리어 거실 침대 탁상 플로어 스탠딩 LED 스탠드 조명 무드등"
$ws.Cells.Item(25, 2).Value = "이케아 HEKTOGRAM 헥토그람 상향식플로어스탠드 조명 거실등 전구포함"
$ws.Cells.Item(26, 2).Value = "마켓비 OGOLD 장스탠드 플로어스탠드 무드등 조명 인테리어"
$ws.Cells.Item(27, 2).Value = "루이스 판텔라 플로어 조명 램프 거실 장 스탠드 유광 고급형 밝기조절 LED 램프2개포함"
$ws.Cells.Item(28, 2).Value = "원목 아폴론 스탠드 조명 장스탠드 무드등 선반 램프"
$ws.Cells.Item(29, 2).Value = "삼성전자 인테리어 스탠드 플로어용 단스탠드 조명 바디 쉐이드 브릭"
$ws.Cells.Item(30, 2).Value = "가온 LED 장스탠드 플로어스탠드 거실 무드등 조명 인테리어"
$ws.Cells.Item(31, 2).Value = "깃털 장스탠드 플로어램프 단스탠드 무드등 인테리어조명"
$ws.Cells.Item(32, 2).Value = "다이슨 라이트사이클 모프 조명 플로어스탠드형"
$ws.Cells.Item(33, 2).Value = "이케아조명등 활장 거실장스탠드전등 식물등 장스텐드조명 스탠딩조명"
$ws.Cells.Item(34, 2).Value = "거실 장스탠드 독서 조명 장스텐드조명"
$ws.Cells.Item(35, 2).Value = "이케아텔뷘 조명등 장스탠드 스탠드조명 등 긴 플로어 스텐드조명 버섯모양스탠드"
$ws.Cells.Item(36, 2).Value = "조명등 장스탠드 활 식물등 플로어 스탠드 조명 뉘모네 라나르프 레르스타 필립스휴 이케아"
$ws.Cells.Item(37, 2).Value = "이케아장스탠드 스탠드조명 조명등 무드등 거실 플로어 인테리어 스텐드 긴스탠드 예쁜스탠드"
$ws.Cells.Item(38, 2).Value = "이케아 ARSTID 오르스티드 플로어스탠드 장스탠드 거실 조명 스탠드 니켈도금 102.840.36"
$ws.Cells.Item(39, 2).Value = "장스탠드 플로어 활 스탠드 거실 스탠드조명 리모컨"
$ws.Cells.Item(40, 2).Value = "파파조명 파파 로델리 장스탠드 인테리어스탠드 무드등 조명"
$ws.Cells.Item(41, 2).Value = "LED 코너 플로어 램프 거실 서 장식 분위기 램프 침실 장식 스탠드 조명 실내 조명"
$ws.Cells.Item(42, 2).Value = "홈즈 오로라 스탠드 조명 간접 무드등"
$ws.Cells.Item(43, 2).Value = "이케아 테르티알 포르소 책상 집게 스탠드 조명 LED 고정 공부 학습용 독서등 작업등"
$ws.Cells.Item(44, 2).Value = "활장스탠드 장스탠드 거실등 스탠드조명 플로어스탠드"
$ws.Cells.Item(45, 2).Value = "덴마크 카프라니 1960 Caprani 플로어 램프 장스탠드 빈티지 플리츠 레트로 조명"
$ws.Cells.Item(46, 2).Value = "이케아 무드등 LED 장스탠드 조명 식물등 독서등 거실 인테리어"
$ws.Cells.Item(47, 2).Value = "마켓비 장스탠드 엔틱 롱 원룸 장스텐드 조명 등 거실"
$ws.Cells.Item(48, 2).Value = "거실 스탠드조명 북유럽 창의적이고 개성 서재 골드 모던 심플"
$ws.Cells.Item(49, 2).Value = "스탠딩조명 거실 안방 침실 소파 장스탠드 클래식 포인트 스탠드 조명"
$ws.Cells.Item(50, 2).Value = "마켓비 장스탠드 침실 거실 조명 램프 롱 활장 플로어 독서등 무드등 인테리어"
$ws.Cells.Item(51, 2).Value = "IKEA 이케아 ÅRSTID 오르스티드 플로어스탠드 장스탠드 거실조명 니켈도금, 황동"
$ws.Cells.Item(52, 2).Value = "북유럽 거실 침실 스탠드조명 심플 인테리어 거실 소파등 북유럽 심플 스탠드 이케아 장식"
$ws.Cells.Item(53, 2).Value = "[오늘 출발] 이케아 플로어 스탠드 조명 인테리어 레르스타 거실 침실 장스탠드 독서등"
$ws.Cells.Item(54, 2).Value = "LED독서등 눈보호 LED스탠드 조명 책상 공부 스탠드 미니"
$ws.Cells.Item(55, 2).Value = "루이스 폴센 판텔라  플로어 고급형 스탠드 조명 램프 장스탠드  무드등 수입 조명"
$ws.Cells.Item(56, 2).Value = "바우하우스 카이저이델 모던 스탠드 조명 침대 옆 테이블 램프 복고풍"
$ws.Cells.Item(57, 2).Value = "스피아노 몬스터램프 거실조명 치히로의섬 키다리스탠드"
$ws.Cells.Item(58, 2).Value = "앤틱 거실 스탠드 조명 등 바로니스 플로어 장스탠드 LED 전등 결혼선물 집들이선물"
$ws.Cells.Item(59, 2).Value = "IKEA 장스탠드 LED 스탠드 조명 독서등 무드등 거실등"
$ws.Cells.Item(60, 2).Value = "이케아 TALLBYN 텔뷘 플로어스탠드/전구미포함/조명"
$ws.Cells.Item(61, 2).Value = "캥거 장스탠드 만달라키 노을 석양 조명 선셋 무드등 MEL001"
$ws.Cells.Item(62, 2).Value = "벨라 장스탠드 거실스탠드 플로어 조명 무드등 인테리어스탠드"
$ws.Cells.Item(63, 2).Value = "책상 침대 무선 집게스탠드 조명 독서등 미니 램프 테이블 북라이트"
$ws.Cells.Item(64, 2).Value = "크리스탈 장스탠드 LED 조명 거실 침실 인테리어 조명등기구"
$ws.Cells.Item(65, 2).Value = "이케아 TOMELILLA 토멜릴라 플로어스탠드 조명 전구미포함"
$ws.Cells.Item(66, 2).Value = "포라이트 플라워팟 VP5 플로어 거실 장 스탠드 박지성조명"
$ws.Cells.Item(67, 2).Value = "루나 장 스탠드 조명"
$ws.Cells.Item(68, 2).Value = "조명 에시아 장스탠드 조명스탠드 LED스탠드"
$ws.Cells.Item(69, 2).Value = "북유럽 모던 스탠드 조명 북유럽 포스트 모던 산업 스"
$ws.Cells.Item(70, 2).Value = "북유럽 스탠드 루이스 판텔라 플로어 램프 장 스탠드 조명 오팔화이트 거실 소파 테이블 모던"
$ws.Cells.Item(71, 2).Value = "스탠드조명 장스탠드 침실조명 트리니티 하프 장스탠드"
$ws.Cells.Item(72, 2).Value = "복고풍 레트로 크리스탈 스탠드 조명 램프 앤틱스타일"
$ws.Cells.Item(73, 2).Value = "LED장스탠드 스탠드조명 장스탠드 코너 스탠드 램프 RGB 3색 리모콘"
$ws.Cells.Item(74, 2).Value = "마켓비 단스탠드 라탄 이케아 조명 책상 미니 독서등 침대 침실 우드 엔틱 무드등 갓 램프"
$ws.Cells.Item(75, 2).Value = "이케아 SIMRISHAMN 심리스함 플로어스탠드 전구미포함 조명 장스탠드 거실등"
$ws.Cells.Item(76, 2).Value = "이자벨32 클래식장스탠드조명 장스텐드조명 스탠드"
$ws.Cells.Item(77, 2).Value = "이케아 LERSTA 레르스타 플로어스탠드 스마트 조명 독서등 장 스탠드 거실 2color"
$ws.Cells.Item(78, 2).Value = "루이스폴센 판텔라 조명 인테리어 장 스탠드 스틸"
$ws.Cells.Item(79, 2).Value = "선셋 만달라키 플로어 스탠드 조명 석양 감성 무드등"
$ws.Cells.Item(80, 2).Value = "거실장스탠드 스탠드조명 스탠드등 이케아 스탠딩 식물 침실 긴 활장 플로어램프 롱 인테리어"
$ws.Cells.Item(81, 2).Value = "인테리어 장스탠드 취침등 수유등 플로어스탠드 롱스탠드 조명"
$ws.Cells.Item(82, 2).Value = "이케아 얀셰, 집게 스탠드조명 클립 자바라 led독서등"
$ws.Cells.Item(83, 2).Value = "스탠딩조명 STRIT 거실장스탠드 식물등스탠드 필립스휴스탠드 E26 골드 무드등 마켓비"
$ws.Cells.Item(84, 2).Value = "스탠드조명 장스탠드 침실조명 코니 아일랜드 장스탠드"
$ws.Cells.Item(85, 2).Value = "장스탠드 플로어 활 스탠드 거실 스탠드조명 리모컨"
$ws.Cells.Item(86, 2).Value = "파파 튤립 장스탠드 인테리어스탠드 조명"
$ws.Cells.Item(87, 2).Value = "거실스탠드조명 플로어 램프 ins 북유럽의 창의적인 아이덴티티_ 000015350"
$ws.Cells.Item(88, 2).Value = "북유럽 거실 플로어 깃털스탠드 인테리어 조명"
$ws.Cells.Item(89, 2).Value = "한샘 한샘몰X스피아노 헤이즈 스탠드 조명SET 플로어 테이블 인테리어스탠드 장스탠드"
$ws.Cells.Item(90, 2).Value = "라루즈 라르고 밝기조절 장스탠드 무드등 식탁 스탠드조명 거실등 인테리어 조명 LED 디밍"
$ws.Cells.Item(91, 2).Value = "스탠드 조명 북유럽 스탠드 거실의 미니멀리즘과 심플한 침실 서재 침대의 인스풍 스탠드"
$ws.Cells.Item(92, 2).Value = "이케아 가성비갑 플로어 장 거실 인테리어 스탠드 독서등 조명 등기구 LED 북유럽 레르스타 오르스티드 라나르프 홀뫼"
$ws.Cells.Item(93, 2).Value = "마켓비 집게 스탠드 단 집게형 조명 LED 독서등 책상 테이블 미니램프"
$ws.Cells.Item(94, 2).Value = "엘사 크리스탈 눈꽃스탠드 조명"
$ws.Cells.Item(95, 2).Value = "Flos Superloon 플로스 디자인 조명 플로어 스탠드 거실 장스탠드 슈퍼룬"
$ws.Cells.Item(96, 2).Value = "북유럽 낚시등 디자인 인테리어 활 장스탠드 스탠드조명"
$ws.Cells.Item(97, 2).Value = "루이스폴센 루이스폴센 LOUIS POULSEN PH 3½-2½ FLOOR lamp 스탠드 조명 chrome"
$ws.Cells.Item(98, 2).Value = "장스탠드 활장스탠드 스탠딩조명 장조명 거실장스탠드"
$ws.Cells.Item(99, 2).Value = "코너 II 침실스탠드 LED조명 인테리어 거실 호텔 카페"
$ws.Cells.Item(100, 2).Value = "LED 플로어 램프 거실 침실 책상 스탠드 조명 눈보호"
$ws.Cells.Item(101, 2).Value = "프레디 장 스탠드 조명 탁상조명 스탠드"

# Rows 43/44 also traded places in the source feed (everything except
# the running index in column A) - swap C:N via Copy so numeric-looking
# text (prices) keeps its original text typing instead of becoming a number.
$ws.Range("C44:N44").Copy($ws.Range("P1:AA1"))
$ws.Range("C43:N43").Copy($ws.Range("C44:N44"))
$ws.Range("P1:AA1").Copy($ws.Range("C43:N43"))
$ws.Range("P1:AA1").ClearContents()

# Copy() of a blank source cell leaves the destination untouched, so
# explicitly blank out the brand columns that came from row 43's blanks.
$ws.Range("I44").ClearContents()
$ws.Range("J44").ClearContents()
